$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns Y:AE
$ws.Range("Y1").Value = "big-screen"
$ws.Range("Z1").Value = "small-screen"
$ws.Range("AA1").Value = "too-heavy"
$ws.Range("AB1").Value = "high-performance"
$ws.Range("AC1").Value = "low-performance"
$ws.Range("AD1").Value = "too-expensive"
$ws.Range("AE1").Value = "to-cheap"

# New data values, rows 2-6
$data = @(
    @(9, 5, 6, 7.5, 6.5, 3, 9),
    @(8, 6, 7, 7, 7, 4, 8),
    @(6, 7, 8, 6.5, 7.5, 6, 6),
    @(5, 8, 7.5, 6, 8, 8, 4),
    @(8, 6, 3, 9, 5, 7, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 25).Value = $vals[0]
    $ws.Cells.Item($row, 26).Value = $vals[1]
    $ws.Cells.Item($row, 27).Value = $vals[2]
    $ws.Cells.Item($row, 28).Value = $vals[3]
    $ws.Cells.Item($row, 29).Value = $vals[4]
    $ws.Cells.Item($row, 30).Value = $vals[5]
    $ws.Cells.Item($row, 31).Value = $vals[6]
}

# Update column widths to best-fit for new text columns (nearest value the
# host's integer-pixel column-width model can represent)
$ws.Range("Z1").ColumnWidth = 12.285714285714286
$ws.Range("AB1").ColumnWidth = 18
$ws.Range("AC1").ColumnWidth = 16.857142857142858
$ws.Range("AD1").ColumnWidth = 14

# Update view: scroll the window right and move the selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 24
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AE14").Select()
